$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move existing content from old columns (C, I) to new columns (B, H) ---
# These strings already exist in the shared string table, so the order in
# which we touch them does not affect the final shared-string ordering.

$ws.Range("C4").ClearContents()
$ws.Range("B4").Value = "показать отчет по депозиту"

$ws.Range("C5").ClearContents()
$ws.Range("B5").Value = "DepositViewModel"
$ws.Range("B5").Font.Bold = $true

$ws.Range("C7").ClearContents()
$ws.Range("B10").Value = "показать сводную форму по всем депозитам"

$ws.Range("C8").ClearContents()
$ws.Range("B11").Value = "DepositsViewModel"
$ws.Range("B11").Font.Bold = $true

$ws.Range("C10").ClearContents()
$ws.Range("B17").Value = "показать ожидаемые доходы от депозитов"

$ws.Range("C11").ClearContents()
$ws.Range("B18").Value = "MonthAnalysisViewModel"
$ws.Range("B18").Font.Bold = $true

$ws.Range("I4").ClearContents()
$ws.Range("H4").Value = "DepositExtractor"
$ws.Range("H4").Font.Bold = $true

$ws.Range("I5").ClearContents()
$ws.Range("H5").Value = "Находит все операции по данному счету"

$ws.Range("I6").ClearContents()
$ws.Range("H6").Value = "и составляет таблицу ежедневных остатков"

$ws.Range("I7").ClearContents()
$ws.Range("H7").Value = "и общие суммы взносов, процентов, расходов"

# --- Add the brand-new cells/strings, in the exact order needed so the
# shared-string table ends up with the new entries appended in this order ---

$ws.Range("B6").Value = "нужна статистика и "
$ws.Range("B7").Value = "нужен прогноз по месяцу и до конца"
$ws.Range("B12").Value = "нужна статистика и"
$ws.Range("B13").Value = "нужно определение какие %%"
$ws.Range("B14").Value = "относятся к какому году"
$ws.Range("B20").Value = "нужен прогноз по месяцу"

$ws.Range("H3").Value = "статистика"

$ws.Range("O5").Value = "отчеты"

$ws.Range("O6").Value = "DepositReporter"
$ws.Range("O6").Font.Bold = $true

$ws.Range("O8").Value = "DepositExcelReporter"
$ws.Range("O8").Font.Bold = $true

$ws.Range("O7").Value = "составляет List<String> для отчета"
$ws.Range("O9").Value = "составляет файл экселя"

$ws.Range("H10").Value = "расчет"

$ws.Range("H11").Value = "DepositCalculator"
$ws.Range("H11").Font.Bold = $true

# --- Reuse of already-created string "нужна статистика и" ---
$ws.Range("B19").Value = "нужна статистика и"

# --- Selection matches the diff's new active cell ---
$ws.Range("H11").Select()

# --- Page setup: paper size / orientation ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
